# Backbone energy assets: add capacity_electric_kw / capacity_heat_kw columns
# to the config_energyAssets sheet.
#
# - existing "capacity_kw" column (G) is renamed to "capacity_electric_kw"
# - a new "capacity_heat_kw" column (H) is appended, defaulting to 0 for the
#   existing energy-asset rows
# - selection/dimension bookkeeping is updated to match

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("config_energyAssets")

# Rename the existing capacity column to be electric-specific ...
$ws.Range("G1").Value = "capacity_electric_kw"

# ... and add the new heat-capacity column next to it.
$ws.Range("H1").Value = "capacity_heat_kw"

# energyAssets rows 2-7 all get a default heat capacity of 0 kW (these are
# electric production assets - windmills/photovoltaic - with no heat output).
$lastRow = 7
for ($row = 2; $row -le $lastRow; $row++) {
    $ws.Cells.Item($row, 8).Value = 0
}

# Move the selection like the author's session did.
$ws.Range("E5").Select() | Out-Null
